# Natmi following Dr Hou advice
# Update Ligand/Receptor-expressing-cells counts (1 -> 3) and recompute
# dependent expression/specificity statistics for rows 2-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 100.8373843333333
$ws.Range("H2").Value = 302.512153
$ws.Range("I2").Value = 0.6551985585448407
$ws.Range("J2").Value = 0.6551985585448408
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.16653
$ws.Range("N2").Value = 51.49959
$ws.Range("O2").Value = 0.0560345397128279
$ws.Range("P2").Value = 0.0560345397128279
$ws.Range("Q2").Value = 1731.027983279696
$ws.Range("R2").Value = 15579.25184951727
$ws.Range("S2").Value = 0.03671374964856847
$ws.Range("T2").Value = 0.03671374964856847

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 100.8373843333333
$ws.Range("H3").Value = 302.512153
$ws.Range("I3").Value = 0.6551985585448407
$ws.Range("J3").Value = 0.6551985585448408
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 256.4443053333333
$ws.Range("N3").Value = 769.332916
$ws.Range("O3").Value = 0.8370788162388805
$ws.Range("P3").Value = 0.8370788162388805
$ws.Range("Q3").Value = 25859.17297699201
$ws.Range("R3").Value = 232732.5567929281
$ws.Range("S3").Value = 0.548452833788136
$ws.Range("T3").Value = 0.5484528337881361

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 100.8373843333333
$ws.Range("H4").Value = 302.512153
$ws.Range("I4").Value = 0.6551985585448407
$ws.Range("J4").Value = 0.6551985585448408
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.74538866666666
$ws.Range("N4").Value = 98.236166
$ws.Range("O4").Value = 0.1068866440482915
$ws.Range("P4").Value = 0.1068866440482915
$ws.Range("Q4").Value = 3301.959342125044
$ws.Range("R4").Value = 29717.6340791254
$ws.Range("S4").Value = 0.07003197510813608
$ws.Range("T4").Value = 0.07003197510813609

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 23.90796933333333
$ws.Range("H5").Value = 71.72390799999999
$ws.Range("I5").Value = 0.1553438454249564
$ws.Range("J5").Value = 0.1553438454249564
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.16653
$ws.Range("N5").Value = 51.49959
$ws.Range("O5").Value = 0.0560345397128279
$ws.Range("P5").Value = 0.0560345397128279
$ws.Range("Q5").Value = 410.4168727997466
$ws.Range("R5").Value = 3693.75185519772
$ws.Range("S5").Value = 0.008704620875608119
$ws.Range("T5").Value = 0.008704620875608119

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 23.90796933333333
$ws.Range("H6").Value = 71.72390799999999
$ws.Range("I6").Value = 0.1553438454249564
$ws.Range("J6").Value = 0.1553438454249564
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 256.4443053333333
$ws.Range("N6").Value = 769.332916
$ws.Range("O6").Value = 0.8370788162388805
$ws.Range("P6").Value = 0.8370788162388805
$ws.Range("Q6").Value = 6131.062587617302
$ws.Range("R6").Value = 55179.56328855572
$ws.Range("S6").Value = 0.1300350422383181
$ws.Range("T6").Value = 0.1300350422383181

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 23.90796933333333
$ws.Range("H7").Value = 71.72390799999999
$ws.Range("I7").Value = 0.1553438454249564
$ws.Range("J7").Value = 0.1553438454249564
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 32.74538866666666
$ws.Range("N7").Value = 98.236166
$ws.Range("O7").Value = 0.1068866440482915
$ws.Range("P7").Value = 0.1068866440482915
$ws.Range("Q7").Value = 782.8757480507473
$ws.Range("R7").Value = 7045.881732456727
$ws.Range("S7").Value = 0.01660418231103014
$ws.Range("T7").Value = 0.01660418231103014

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 29.15819666666667
$ws.Range("H8").Value = 87.47459000000001
$ws.Range("I8").Value = 0.1894575960302029
$ws.Range("J8").Value = 0.1894575960302029
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 17.16653
$ws.Range("N8").Value = 51.49959
$ws.Range("O8").Value = 0.0560345397128279
$ws.Range("P8").Value = 0.0560345397128279
$ws.Range("Q8").Value = 500.5450578242333
$ws.Range("R8").Value = 4504.9055204181
$ws.Range("S8").Value = 0.01061616918865131
$ws.Range("T8").Value = 0.01061616918865131

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 29.15819666666667
$ws.Range("H9").Value = 87.47459000000001
$ws.Range("I9").Value = 0.1894575960302029
$ws.Range("J9").Value = 0.1894575960302029
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 256.4443053333333
$ws.Range("N9").Value = 769.332916
$ws.Range("O9").Value = 0.8370788162388805
$ws.Range("P9").Value = 0.8370788162388805
$ws.Range("Q9").Value = 7477.453488956049
$ws.Range("R9").Value = 67297.08140060444
$ws.Range("S9").Value = 0.1585909402124263
$ws.Range("T9").Value = 0.1585909402124263

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 29.15819666666667
$ws.Range("H10").Value = 87.47459000000001
$ws.Range("I10").Value = 0.1894575960302029
$ws.Range("J10").Value = 0.1894575960302029
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 32.74538866666666
$ws.Range("N10").Value = 98.236166
$ws.Range("O10").Value = 0.1068866440482915
$ws.Range("P10").Value = 0.1068866440482915
$ws.Range("Q10").Value = 954.7964826691044
$ws.Range("R10").Value = 8593.16834402194
$ws.Range("S10").Value = 0.02025048662912531
$ws.Range("T10").Value = 0.02025048662912531
